$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 143, shifting existing rows 143:163 down to 144:164
$ws.Rows("143:143").Insert()

# Populate the newly inserted row 143 with its data
$ws.Range("A143").Value = "SOLEMON2025"
$ws.Range("B143").Value = "ITA17"
$ws.Range("C143").Value = "35"
$ws.Range("D143").Value = "2-RAP"
$ws.Range("E143").Value = "PAPELON"
$ws.Range("F143").Value = 1
$ws.Range("G143").Value = 1
$ws.Range("H143").Value = 19
$ws.Range("I143").Value = 5
$ws.Range("J143").Value = "F"
$ws.Range("K143").Value = "MEDPF-1"
$ws.Range("L143").Value = 0
$ws.Range("P143").Value = "NA l inferred"
